# Updated cryptos list values (Price and Volume(1h) columns) to reflect
# the latest scrape. Also corrects the ordering of the Quant / Cronos rows
# (rows 49 and 50 swap which coin they describe, matching the fetched API order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.366.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.81%  '
$ws.Range("D3").Value = "'1.819.33"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.94%  '
$ws.Range("D5").Value = "'330.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.88%  '
$ws.Range("E6").Value = '  -0.70%  '
$ws.Range("D7").Value = "'0.4568"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.59%  '
$ws.Range("E8").Value = '  -4.27%  '
$ws.Range("D9").Value = "'46.11"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.27%  '
$ws.Range("D11").Value = "'0.9610"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'20.95"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.95%  '
$ws.Range("D13").Value = "'1.848.35"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.92%  '
$ws.Range("D14").Value = "'5.852"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.84%  '
$ws.Range("D15").Value = "'7.066"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.27%  '
$ws.Range("D16").Value = "'1.003"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.95%  '
$ws.Range("D17").Value = "'89.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.58%  '
$ws.Range("D18").Value = "'0.06587"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.23%  '
$ws.Range("D19").Value = "'0.00001020"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.83%  '
$ws.Range("D20").Value = "'17.09"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.89%  '
$ws.Range("D21").Value = "'1.002"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("D22").Value = "'27.361.07"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.92%  '
$ws.Range("D23").Value = "'5.294"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.00%  '
$ws.Range("E24").Value = '  -2.03%  '
$ws.Range("D25").Value = "'2.271"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.99%  '
$ws.Range("D26").Value = "'2.052.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.74%  '
$ws.Range("D27").Value = "'155.94"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.17%  '
$ws.Range("D28").Value = "'19.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.97%  '
$ws.Range("D29").Value = "'2.040"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -6.06%  '
$ws.Range("D30").Value = "'5.234"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.24%  '
$ws.Range("D31").Value = "'117.66"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.76%  '
$ws.Range("D32").Value = "'0.09308"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("D33").Value = "'0.9320"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.75%  '
$ws.Range("D34").Value = "'3.569"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.91%  '
$ws.Range("D35").Value = "'5.212"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.67%  '
$ws.Range("D36").Value = "'1.307"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.13%  '
$ws.Range("D37").Value = "'0.05913"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.01%  '
$ws.Range("D38").Value = "'0.02169"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.69%  '
$ws.Range("D39").Value = "'8.077"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.72%  '
$ws.Range("D40").Value = "'1.002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("E41").Value = '  -5.61%  '
$ws.Range("D42").Value = "'0.5747"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.39%  '
$ws.Range("D43").Value = "'0.1814"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.56%  '
$ws.Range("D44").Value = "'9.885"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.40%  '
$ws.Range("D45").Value = "'1.279"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("D46").Value = "'11.90"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.62%  '
$ws.Range("D47").Value = "'0.5395"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.12%  '
$ws.Range("D48").Value = "'1.866"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -4.37%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'0.06556"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.60%  '
$ws.Range("B50").Value = 'Quant'
$ws.Range("C50").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D50").Value = "'109.73"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("E51").Value = '  -34.24%  '
